$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Flip the Execution Flag column (C) values: most rows NO -> YES,
# row 19 YES -> NO. Rows 7, 11, 14-18 remain NO (unchanged).
$ws.Range("C2").Value = "YES"
$ws.Range("C3").Value = "YES"
$ws.Range("C4").Value = "YES"
$ws.Range("C5").Value = "YES"
$ws.Range("C6").Value = "YES"
$ws.Range("C8").Value = "YES"
$ws.Range("C9").Value = "YES"
$ws.Range("C10").Value = "YES"
$ws.Range("C12").Value = "YES"
$ws.Range("C13").Value = "YES"
$ws.Range("C19").Value = "NO"
